# Updated cryptos list on Fri Feb 24 13:41:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text in the
# source data (e.g. "23.938.10", "1.000"). Excel's COM layer auto-coerces a
# bare numeric-looking string assigned to .Value into a real number, which
# would corrupt values like "1.004" (-> 1.004, fine) but especially
# multi-dot figures like "23.808.65" (-> parsed oddly) and strips
# significant trailing zeros ("0.3880" -> 0.388). Force the target cells to
# Text format first so the assigned strings are preserved verbatim.
$priceCells = @(
    "D2","D3","D4","D5","D6","D7","D8","D9","D10",
    "D11","D12","D13","D14","D15","D16","D17","D18","D19","D20",
    "D21","D23","D24","D25","D26","D27","D28","D29","D30",
    "D31","D32","D33","D34","D35","D36","D37","D38","D39","D40",
    "D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2 : Bitcoin ---
$ws.Range("D2").Value = "23.808.65"
$ws.Range("E2").Value = "  -1.30%  "

# --- Row 3 : Ethereum ---
$ws.Range("D3").Value = "1.641.11"
$ws.Range("E3").Value = "  -1.24%  "

# --- Row 4 : TetherUSD ---
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.60%  "

# --- Row 5 : BNB ---
$ws.Range("D5").Value = "309.74"
$ws.Range("E5").Value = "  -0.10%  "

# --- Row 6 : USDC ---
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.44%  "

# --- Row 7 : XRP ---
$ws.Range("D7").Value = "0.3880"
$ws.Range("E7").Value = "  -1.12%  "

# --- Row 8 : Cardano ---
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").Value = "  -1.95%  "

# --- Row 9 : OKB ---
$ws.Range("D9").Value = "50.47"
$ws.Range("E9").Value = "  -2.70%  "

# --- Row 10 : Polygon ---
$ws.Range("D10").Value = "1.323"
$ws.Range("E10").Value = "  -4.20%  "

# --- Row 11 : BinanceUSD ---
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +0.36%  "

# --- Row 12 : Dogecoin ---
$ws.Range("D12").Value = "0.08386"
$ws.Range("E12").Value = "  -1.47%  "

# --- Row 13 : Solana ---
$ws.Range("D13").Value = "23.66"
$ws.Range("E13").Value = "  -2.69%  "

# --- Row 14 : Polkadot ---
$ws.Range("D14").Value = "6.954"
$ws.Range("E14").Value = "  -4.60%  "

# --- Row 15 : Chainlink ---
$ws.Range("D15").Value = "7.834"
$ws.Range("E15").Value = "  -4.11%  "

# --- Row 16 : ShibaInu ---
$ws.Range("D16").Value = "0.00001308"
$ws.Range("E16").Value = "  -1.20%  "

# --- Row 17 : WrappedEther ---
$ws.Range("D17").Value = "1.639.04"
$ws.Range("E17").Value = "  -0.94%  "

# --- Row 18 : Litecoin ---
$ws.Range("D18").Value = "93.52"
$ws.Range("E18").Value = "  -1.74%  "

# --- Row 19 : TRON ---
$ws.Range("D19").Value = "0.06960"
$ws.Range("E19").Value = "  -0.27%  "

# --- Row 20 : Avalanche ---
$ws.Range("D20").Value = "19.40"
$ws.Range("E20").Value = "  -3.65%  "

# --- Row 21 : Uniswap ---
$ws.Range("D21").Value = "6.869"
$ws.Range("E21").Value = "  -1.76%  "

# --- Row 22 : Dai (price unchanged, volume changed) ---
$ws.Range("E22").Value = "  +0.18%  "

# --- Row 23 : Cosmos ---
$ws.Range("D23").Value = "13.56"
$ws.Range("E23").Value = "  -1.48%  "

# --- Row 24 : WrappedBTC ---
$ws.Range("D24").Value = "23.819.61"
$ws.Range("E24").Value = "  -1.22%  "

# --- Row 25 : Toncoin ---
$ws.Range("D25").Value = "2.437"
$ws.Range("E25").Value = "  -3.03%  "

# --- Row 26 : LidoDAOToken ---
$ws.Range("D26").Value = "2.872"
$ws.Range("E26").Value = "  -9.68%  "

# --- Row 27 : EthereumClassic ---
$ws.Range("D27").Value = "21.86"
$ws.Range("E27").Value = "  -2.28%  "

# --- Rows 28/29 : Monero <-> HuobiToken swap positions ---
$ws.Range("B28").Value = "HuobiToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D28").Value = "5.587"
$ws.Range("E28").Value = "  +4.87%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "153.04"
$ws.Range("E29").Value = "  -0.21%  "

# --- Row 30 : BitcoinCash ---
$ws.Range("D30").Value = "136.54"
$ws.Range("E30").Value = "  -2.41%  "

# --- Row 31 : WEMIXTOKEN ---
$ws.Range("D31").Value = "2.496"
$ws.Range("E31").Value = "  -0.01%  "

# --- Row 32 : Filecoin ---
$ws.Range("D32").Value = "7.611"
$ws.Range("E32").Value = "  -3.84%  "

# --- Row 33 : WrappedliquidstakedEther2.0 ---
$ws.Range("D33").Value = "1.823.55"
$ws.Range("E33").Value = "  -0.70%  "

# --- Row 34 : Hedera ---
$ws.Range("D34").Value = "0.07979"
$ws.Range("E34").Value = "  -2.35%  "

# --- Row 35 : ImmutableX ---
$ws.Range("D35").Value = "0.9746"
$ws.Range("E35").Value = "  -7.58%  "

# --- Row 36 : VeChain ---
$ws.Range("D36").Value = "0.02890"
$ws.Range("E36").Value = "  -4.91%  "

# --- Row 37 : InternetComputer(DFINITY) ---
$ws.Range("D37").Value = "6.546"
$ws.Range("E37").Value = "  -3.31%  "

# --- Row 38 : Algorand ---
$ws.Range("D38").Value = "0.2652"
$ws.Range("E38").Value = "  -3.31%  "

# --- Rows 39/40 : FraxShare <-> Stellar swap positions ---
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.09078"
$ws.Range("E39").Value = "  -1.34%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "10.32"
$ws.Range("E40").Value = "  -8.32%  "

# --- Row 41 : TheSandbox ---
$ws.Range("D41").Value = "0.7478"
$ws.Range("E41").Value = "  -2.25%  "

# --- Row 42 : Aptos ---
$ws.Range("D42").Value = "13.25"
$ws.Range("E42").Value = "  -2.45%  "

# --- Row 43 : TrustWalletToken ---
$ws.Range("D43").Value = "1.414"
$ws.Range("E43").Value = "  -0.92%  "

# --- Row 44 : EnergySwap ---
$ws.Range("D44").Value = "16.45"
$ws.Range("E44").Value = "  -1.32%  "

# --- Row 45 : Decentraland ---
$ws.Range("D45").Value = "0.6876"
$ws.Range("E45").Value = "  -2.58%  "

# --- Row 46 : NEARProtocol ---
$ws.Range("D46").Value = "2.414"
$ws.Range("E46").Value = "  -4.26%  "

# --- Row 47 : PancakeSwap ---
$ws.Range("D47").Value = "4.082"
$ws.Range("E47").Value = "  -0.32%  "

# --- Row 48 : Frax ---
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.42%  "

# --- Row 49 : Cronos ---
$ws.Range("D49").Value = "0.08214"
$ws.Range("E49").Value = "  -2.19%  "

# --- Row 50 : Quant ---
$ws.Range("D50").Value = "133.87"
$ws.Range("E50").Value = "  -1.60%  "

# --- Row 51 : Flow ---
$ws.Range("D51").Value = "1.211"
$ws.Range("E51").Value = "  -3.14%  "
